$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Add two new files (3f14bde9-... and ba0dd1f4-...) to the handback
# report: one new row per file on "Overview", "zh-cn" and "de-de".
# ---------------------------------------------------------------------

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$loOverview = $overview.ListObjects.Item(1)
$loZhcn     = $zhcn.ListObjects.Item(1)
$loDede     = $dede.ListObjects.Item(1)

# Grow each table by two rows (keeps table ref / autofilter / dimension
# in sync, same as Excel does when a ListRow is appended).
$loOverview.ListRows.Add() | Out-Null
$loOverview.ListRows.Add() | Out-Null
$loZhcn.ListRows.Add() | Out-Null
$loZhcn.ListRows.Add() | Out-Null
$loDede.ListRows.Add() | Out-Null
$loDede.ListRows.Add() | Out-Null

# --------------------------- Overview sheet ---------------------------
# Columns: A=File Name, B=Path And Name, C=Extension, D=Publish URL,
#          E=zh-cn, F=de-de, G=Latest HO Xliff Generate Date
$overview.Cells.Item(4, 1).Value = "3f14bde9-0c56-4a70-af96-9a02074d8301.md"
$overview.Cells.Item(4, 2).Value = "e2e\3f14bde9-0c56-4a70-af96-9a02074d8301.md"
$overview.Cells.Item(4, 3).Value = ".md"
$overview.Cells.Item(4, 4).Value = ""
$overview.Cells.Item(4, 5).Value = "Ready for handoff"
$overview.Cells.Item(4, 6).Value = "Ready for handoff"
$overview.Cells.Item(4, 7).Value = "2016-08-29 22:42:50"

$overview.Cells.Item(5, 1).Value = "ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md"
$overview.Cells.Item(5, 2).Value = "e2e\ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md"
$overview.Cells.Item(5, 3).Value = ".md"
$overview.Cells.Item(5, 4).Value = ""
$overview.Cells.Item(5, 5).Value = "Ready for handoff"
$overview.Cells.Item(5, 6).Value = "Ready for handoff"
$overview.Cells.Item(5, 7).Value = "2016-08-29 22:42:50"

$overview.Hyperlinks.Add($overview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0d2737df5b05f8f9622cdb5ffc18a7ea43b5464b/e2e/3f14bde9-0c56-4a70-af96-9a02074d8301.md", "", "", "e2e\3f14bde9-0c56-4a70-af96-9a02074d8301.md") | Out-Null
$overview.Hyperlinks.Add($overview.Range("B5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/423cd28fb57506b2649e5eb8b503152ad7233acc/e2e/ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md", "", "", "e2e\ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md") | Out-Null

# ----------------------------- zh-cn sheet -----------------------------
# Columns: A=Source File Name, B=File Extension, C=Status, D=Source Path,
#          E=Priority, F=Content Duplicate, G=Latest Handoff File,
#          H=Latest Handoff Datetime, I=Latest Target File,
#          J=Latest Handback File, K=Latest Handback DateTime,
#          L=Reference Tokens, M=To be localized, N=Dependency From,
#          O=Has metadata, P=Error Detail
$zhcn.Cells.Item(4, 1).Value  = "3f14bde9-0c56-4a70-af96-9a02074d8301.md"
$zhcn.Cells.Item(4, 2).Value  = ".md"
$zhcn.Cells.Item(4, 3).Value  = "Ready for handoff"
$zhcn.Cells.Item(4, 4).Value  = "e2e"
$zhcn.Cells.Item(4, 5).Value  = "ht"
$zhcn.Cells.Item(4, 6).Value  = "False"
$zhcn.Cells.Item(4, 7).Value  = "3f14bde9-0c56-4a70-af96-9a02074d8301.0d2737df5b05f8f9622cdb5ffc18a7ea43b5464b.zh-cn.xlf"
$zhcn.Cells.Item(4, 8).Value  = "2016-08-29 22:42:45"
$zhcn.Cells.Item(4, 9).Value  = ""
$zhcn.Cells.Item(4, 10).Value = ""
$zhcn.Cells.Item(4, 11).Value = "0001-01-01 00:00:00"
$zhcn.Cells.Item(4, 12).Value = ""
$zhcn.Cells.Item(4, 13).Value = "True"
$zhcn.Cells.Item(4, 14).Value = ""
$zhcn.Cells.Item(4, 15).Value = "False"
$zhcn.Cells.Item(4, 16).Value = ""

$zhcn.Cells.Item(5, 1).Value  = "ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md"
$zhcn.Cells.Item(5, 2).Value  = ".md"
$zhcn.Cells.Item(5, 3).Value  = "Ready for handoff"
$zhcn.Cells.Item(5, 4).Value  = "e2e"
$zhcn.Cells.Item(5, 5).Value  = "ht"
$zhcn.Cells.Item(5, 6).Value  = "False"
$zhcn.Cells.Item(5, 7).Value  = "ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.423cd28fb57506b2649e5eb8b503152ad7233acc.zh-cn.xlf"
$zhcn.Cells.Item(5, 8).Value  = "2016-08-29 22:42:45"
$zhcn.Cells.Item(5, 9).Value  = ""
$zhcn.Cells.Item(5, 10).Value = ""
$zhcn.Cells.Item(5, 11).Value = "0001-01-01 00:00:00"
$zhcn.Cells.Item(5, 12).Value = ""
$zhcn.Cells.Item(5, 13).Value = "True"
$zhcn.Cells.Item(5, 14).Value = ""
$zhcn.Cells.Item(5, 15).Value = "False"
$zhcn.Cells.Item(5, 16).Value = ""

$zhcn.Hyperlinks.Add($zhcn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0d2737df5b05f8f9622cdb5ffc18a7ea43b5464b/e2e/3f14bde9-0c56-4a70-af96-9a02074d8301.md", "", "", "3f14bde9-0c56-4a70-af96-9a02074d8301.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/423cd28fb57506b2649e5eb8b503152ad7233acc/e2e/ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md", "", "", "ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md") | Out-Null

# ----------------------------- de-de sheet -----------------------------
$dede.Cells.Item(4, 1).Value  = "3f14bde9-0c56-4a70-af96-9a02074d8301.md"
$dede.Cells.Item(4, 2).Value  = ".md"
$dede.Cells.Item(4, 3).Value  = "Ready for handoff"
$dede.Cells.Item(4, 4).Value  = "e2e"
$dede.Cells.Item(4, 5).Value  = "ht"
$dede.Cells.Item(4, 6).Value  = "False"
$dede.Cells.Item(4, 7).Value  = "3f14bde9-0c56-4a70-af96-9a02074d8301.0d2737df5b05f8f9622cdb5ffc18a7ea43b5464b.de-de.xlf"
$dede.Cells.Item(4, 8).Value  = "2016-08-29 22:42:50"
$dede.Cells.Item(4, 9).Value  = ""
$dede.Cells.Item(4, 10).Value = ""
$dede.Cells.Item(4, 11).Value = "0001-01-01 00:00:00"
$dede.Cells.Item(4, 12).Value = ""
$dede.Cells.Item(4, 13).Value = "True"
$dede.Cells.Item(4, 14).Value = ""
$dede.Cells.Item(4, 15).Value = "False"
$dede.Cells.Item(4, 16).Value = ""

$dede.Cells.Item(5, 1).Value  = "ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md"
$dede.Cells.Item(5, 2).Value  = ".md"
$dede.Cells.Item(5, 3).Value  = "Ready for handoff"
$dede.Cells.Item(5, 4).Value  = "e2e"
$dede.Cells.Item(5, 5).Value  = "ht"
$dede.Cells.Item(5, 6).Value  = "False"
$dede.Cells.Item(5, 7).Value  = "ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.423cd28fb57506b2649e5eb8b503152ad7233acc.de-de.xlf"
$dede.Cells.Item(5, 8).Value  = "2016-08-29 22:42:50"
$dede.Cells.Item(5, 9).Value  = ""
$dede.Cells.Item(5, 10).Value = ""
$dede.Cells.Item(5, 11).Value = "0001-01-01 00:00:00"
$dede.Cells.Item(5, 12).Value = ""
$dede.Cells.Item(5, 13).Value = "True"
$dede.Cells.Item(5, 14).Value = ""
$dede.Cells.Item(5, 15).Value = "False"
$dede.Cells.Item(5, 16).Value = ""

$dede.Hyperlinks.Add($dede.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0d2737df5b05f8f9622cdb5ffc18a7ea43b5464b/e2e/3f14bde9-0c56-4a70-af96-9a02074d8301.md", "", "", "3f14bde9-0c56-4a70-af96-9a02074d8301.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("A5"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/423cd28fb57506b2649e5eb8b503152ad7233acc/e2e/ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md", "", "", "ba0dd1f4-e8cf-4703-a553-938aa98d2ae9.md") | Out-Null
